$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates (row 2 / row 3) ---------------------------------
# Write the new "A3" / "D2" values first so new shared-string entries are
# appended in the same order the target workbook uses them:
#   5 = Iphone6TPU, 6 = Iphone8TPU, 7 = Visu/Silicone/Iphone6/*
$ws.Range("A2").Value = "Iphone6TPU"
$ws.Range("A3").Value = "Iphone8TPU"

$ws.Range("B2").Value = 518
$ws.Range("C2").Value = 238
$ws.Range("D2").Value = "Visu/Silicone/Iphone6/*"

$ws.Range("B3").Value = 538
$ws.Range("C3").Value = 238
$ws.Range("D3").Value = "Visu/Silicone/Iphone8/*"

# --- Column width (column A: 13.85546875 -> 18) ---------------------------
# Excel's stored <col width=.../> differs from the COM ColumnWidth by a
# constant ~0.8333 (5/6) padding offset, so back that out to land on 18.
$ws.Columns.Item(1).ColumnWidth = 18 - 5/6

# --- Selection / scroll position ------------------------------------------
# Moves the active cell to C4 and (since it's back at the sheet's natural
# top-left) drops the explicit topLeftCell="C1" scroll anchor.
$ws.Range("C4").Select()
